$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 5.93
$ws.Range("B10").Value = 6.851999999999999
$ws.Range("B12").Value = 6.549000000000001
$ws.Range("D13").Value = -7.534999999999999
$ws.Range("B18").Value = 6.548999999999999
